# Scheduled market-data refresh: update leve price/profit columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM sheets with the latest snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H87").Value = 37849.57
$ws.Range("J87").Value = 37849.57
$ws.Range("L87").Value = 37849.57
$ws.Range("N87").Value = -40345.57
$ws.Range("H90").Value = 37849.57
$ws.Range("J90").Value = 37849.57
$ws.Range("L90").Value = 113548.71
$ws.Range("N90").Value = -126028.71
$ws.Range("H125").Value = 8476.666999999999
$ws.Range("I125").Value = 8219.25
$ws.Range("J125").Value = 8991.5
$ws.Range("K125").Value = 73973.25
$ws.Range("L125").Value = 80923.5
$ws.Range("M125").Value = -71513.25
$ws.Range("N125").Value = -85843.5
$ws.Range("H132").Value = 1168.0834
$ws.Range("I132").Value = 1168.0834
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3504.2502
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -974.2501999999999
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 5173.72
$ws.Range("J138").Value = 6239.8975
$ws.Range("L138").Value = 18719.6925
$ws.Range("N138").Value = -28999.6925
$ws.Range("H141").Value = 8330
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4323
$ws.Range("I32").Value = 4323
$ws.Range("K32").Value = 4323
$ws.Range("M32").Value = -4036
$ws.Range("H37").Value = 25969
$ws.Range("J37").Value = 25969
$ws.Range("L37").Value = 25969
$ws.Range("N37").Value = -26515
$ws.Range("H45").Value = 1724.5
$ws.Range("I45").Value = 1832.6666
$ws.Range("J45").Value = 1400
$ws.Range("K45").Value = 1832.6666
$ws.Range("L45").Value = 1400
$ws.Range("M45").Value = -1455.6666
$ws.Range("N45").Value = -2154
$ws.Range("H74").Value = 2833.9092
$ws.Range("I74").Value = 655.1667
$ws.Range("J74").Value = 5448.4
$ws.Range("K74").Value = 655.1667
$ws.Range("L74").Value = 5448.4
$ws.Range("M74").Value = 218.8333
$ws.Range("N74").Value = -7196.4
$ws.Range("H77").Value = 2833.9092
$ws.Range("I77").Value = 655.1667
$ws.Range("J77").Value = 5448.4
$ws.Range("K77").Value = 3275.8335
$ws.Range("L77").Value = 27242
$ws.Range("M77").Value = 1092.1665
$ws.Range("N77").Value = -35978
$ws.Range("H122").Value = 2721.625
$ws.Range("I122").Value = 2682.9048
$ws.Range("J122").Value = 2992.6667
$ws.Range("K122").Value = 8048.714399999999
$ws.Range("L122").Value = 8978.000100000001
$ws.Range("M122").Value = -5598.714399999999
$ws.Range("N122").Value = -13878.0001
$ws.Range("H132").Value = 4857.9
$ws.Range("I132").Value = 4516.2
$ws.Range("J132").Value = 5199.6
$ws.Range("K132").Value = 13548.6
$ws.Range("L132").Value = 15598.8
$ws.Range("M132").Value = -11018.6
$ws.Range("N132").Value = -20658.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2786.8125
$ws.Range("I20").Value = 2842.9092
$ws.Range("K20").Value = 2842.9092
$ws.Range("M20").Value = -2595.9092
$ws.Range("H134").Value = 1842.2632
$ws.Range("I134").Value = 1685.4706
$ws.Range("J134").Value = 3175
$ws.Range("K134").Value = 5056.4118
$ws.Range("L134").Value = 9525
$ws.Range("M134").Value = -2521.4118
$ws.Range("N134").Value = -14595

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 941.25
$ws.Range("I58").Value = 847.6667
$ws.Range("J58").Value = 2345
$ws.Range("K58").Value = 847.6667
$ws.Range("L58").Value = 2345
$ws.Range("M58").Value = -644.6667
$ws.Range("N58").Value = -2751
$ws.Range("H122").Value = 1746.5
$ws.Range("I122").Value = 1746.5
$ws.Range("K122").Value = 5239.5
$ws.Range("M122").Value = -2789.5
$ws.Range("H132").Value = 3068.7646
$ws.Range("I132").Value = 2452.4
$ws.Range("J132").Value = 3949.2856
$ws.Range("K132").Value = 7357.200000000001
$ws.Range("L132").Value = 11847.8568
$ws.Range("M132").Value = -4827.200000000001
$ws.Range("N132").Value = -16907.8568
$ws.Range("H136").Value = 941.25
$ws.Range("I136").Value = 847.6667
$ws.Range("J136").Value = 2345
$ws.Range("K136").Value = 2543.0001
$ws.Range("L136").Value = 7035
$ws.Range("M136").Value = 6.999899999999798
$ws.Range("N136").Value = -12135

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 352.18182
$ws.Range("I7").Value = 352.66666
$ws.Range("J7").Value = 350
$ws.Range("K7").Value = 1057.99998
$ws.Range("L7").Value = 1050
$ws.Range("M7").Value = -945.9999800000001
$ws.Range("N7").Value = -1274
$ws.Range("H38").Value = 365.66666
$ws.Range("I38").Value = 334.33334
$ws.Range("J38").Value = 397
$ws.Range("K38").Value = 1003.00002
$ws.Range("L38").Value = 1191
$ws.Range("M38").Value = -656.0000200000001
$ws.Range("N38").Value = -1885
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H92").Value = 1390
$ws.Range("J92").Value = 1390
$ws.Range("L92").Value = 4170
$ws.Range("N92").Value = -6666
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H123").Value = 6000
$ws.Range("J123").Value = 6000
$ws.Range("L123").Value = 18000
$ws.Range("N123").Value = -22900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6338.9
$ws.Range("I80").Value = 4477.2
$ws.Range("K80").Value = 4477.2
$ws.Range("M80").Value = -3479.2
$ws.Range("H83").Value = 6338.9
$ws.Range("I83").Value = 4477.2
$ws.Range("K83").Value = 22386
$ws.Range("M83").Value = -17394
$ws.Range("H102").Value = 3721.6667
$ws.Range("I102").Value = 3721.6667
$ws.Range("K102").Value = 3721.6667
$ws.Range("M102").Value = -2099.6667
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
